# Update "horarios" workbook with freshly scraped schedule data.
# Línea 141 - 810 update.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "LP1912": new scrape timestamp 01:10:33, one new arrival row added.
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 01:10:33"
$ws1.Range("A3").Value = "Total filas: 3"

$ws1.Cells.Item(6, 1).Value = "01:10:32"
$ws1.Cells.Item(6, 4).Value = 2

$ws1.Cells.Item(7, 1).Value = "01:10:32"
$ws1.Cells.Item(7, 4).Value = 48

$ws1.Cells.Item(8, 1).Value = "01:10:32"
$ws1.Cells.Item(8, 2).Value = "02:58"
$ws1.Cells.Item(8, 3).Value = "215_ALUAR"
$ws1.Cells.Item(8, 4).Value = 108
$ws1.Cells.Item(8, 5).Value = "LP1912"

# ---------------------------------------------------------------
# Sheet "LP1912-215": new scrape timestamp 01:10:33, one new arrival row added.
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 01:10:33"
$ws2.Range("A3").Value = "Total filas: 2"

$ws2.Cells.Item(6, 1).Value = "01:10:32"
$ws2.Cells.Item(6, 4).Value = 2

$ws2.Cells.Item(7, 1).Value = "01:10:32"
$ws2.Cells.Item(7, 2).Value = "02:58"
$ws2.Cells.Item(7, 3).Value = "215_ALUAR"
$ws2.Cells.Item(7, 4).Value = 108
$ws2.Cells.Item(7, 5).Value = "LP1912"

# ---------------------------------------------------------------
# Sheet "6203-6173": only the scrape timestamp advances.
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 01:10:33"
